$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
  @{Row=2; D="36.538.28"; E="  +0.42%  "},
  @{Row=3; D="1.943.21"; E="  -0.49%  "},
  @{Row=4; E="  -0.08%  "},
  @{Row=5; D="243.05"; E="  +0.04%  "},
  @{Row=6; D="0.612"; E="  -0.24%  "},
  @{Row=7; E="  -0.03%  "},
  @{Row=8; D="57.42"; E="  -0.89%  "},
  @{Row=9; D="0.364"; E="  -1.61%  "},
  @{Row=10; D="0.0809"; E="  -1.57%  "},
  @{Row=11; E="  -0.92%  "},
  @{Row=12; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="21.81"; E="  +0.83%  "},
  @{Row=13; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.226.18"; E="  -0.33%  "},
  @{Row=14; D="0.807"; E="  -2.85%  "},
  @{Row=15; D="13.34"; E="  -1.38%  "},
  @{Row=16; D="5.18"; E="  -2.31%  "},
  @{Row=17; D="1.944.77"; E="  -1.37%  "},
  @{Row=18; D="36.444.23"},
  @{Row=19; D="69.27"; E="  -0.82%  "},
  @{Row=20; D="0.0₃0855"; E="  -1.90%  "},
  @{Row=21; D="227.24"; E="  -1.07%  "},
  @{Row=22; D="4.96"; E="  -1.10%  "},
  @{Row=23; E="  -0.16%  "},
  @{Row=24; D="2.39"; E="  -4.73%  "},
  @{Row=25; D="2.30"; E="  +0.95%  "},
  @{Row=26; D="9.18"; E="  -3.20%  "},
  @{Row=27; D="160.45"; E="  -1.91%  "},
  @{Row=28; D="0.134"; E="  +12.69%  "},
  @{Row=29; E="  -2.07%  "},
  @{Row=30; E="  -0.06%  "},
  @{Row=31; D="4.64"; E="  -1.72%  "},
  @{Row=32; D="1.09"; E="  -5.35%  "},
  @{Row=33; D="0.0616"; E="  -3.03%  "},
  @{Row=34; D="4.17"; E="  -3.59%  "},
  @{Row=35; B="THORChain"; C="https://coinranking.com/coin/ybmU-kKU+thorchain-rune"; D="6.14"; E="  +0.81%  "},
  @{Row=36; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.00"; E="  -0.03%  "},
  @{Row=37; E="  -1.08%  "},
  @{Row=38; E="  +1.93%  "},
  @{Row=39; E="  +12.66%  "},
  @{Row=40; D="0.0988"; E="  +0.63%  "},
  @{Row=41; E="  +1.05%  "},
  @{Row=42; D="0.0209"; E="  -0.54%  "},
  @{Row=43; D="1.15"; E="  -2.97%  "},
  @{Row=44; D="15.76"; E="  +0.46%  "},
  @{Row=45; D="1.341.54"; E="  -0.14%  "},
  @{Row=46; D="1.03"; E="  -1.11%  "},
  @{Row=47; D="86.25"; E="  -2.54%  "},
  @{Row=48; D="7.12"; E="  -3.16%  "},
  @{Row=49; D="2.82"; E="  -0.24%  "},
  @{Row=50; D="2.118.37"; E="  -0.42%  "},
  @{Row=51; D="43.21"; E="  -6.24%  "}
)

foreach ($item in $changes) {
  $r = $item.Row
  foreach ($col in @("B", "C", "D", "E")) {
    if ($item.ContainsKey($col)) {
      $cell = $ws.Range("$col$r")
      $cell.NumberFormat = "@"
      $cell.Value = $item[$col]
      $cell.ClearFormats()
    }
  }
}

Write-Output "Applied $($changes.Count) row updates"